$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$d.Content.Find.Execute("2026-01-16 Friday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2026-01-17 Saturday", 2)

# Update the division problems in the table, cell by cell, using positional
# addressing (Table row/column) and direct Range.Text assignment. Several
# problem strings are duplicated elsewhere in the document (e.g. "92÷2="
# and "32÷3=" both occur twice, once as an old value and once as a new
# value), so Find.Execute (which is not reliably scoped to a cell's Range
# in this runtime) cannot be used safely; writing the cell's Range.Text
# directly guarantees only the intended cell is touched.
$table = $d.Tables.Item(1)

$changes = @(
    @{ Row = 1;  Col = 1; New = "83÷6=" },
    @{ Row = 1;  Col = 2; New = "31÷3=" },
    @{ Row = 1;  Col = 3; New = "16÷7=" },
    @{ Row = 1;  Col = 4; New = "88÷4=" },
    @{ Row = 1;  Col = 5; New = "54÷4=" },

    @{ Row = 5;  Col = 1; New = "55÷4=" },
    @{ Row = 5;  Col = 2; New = "33÷4=" },
    @{ Row = 5;  Col = 3; New = "34÷6=" },
    @{ Row = 5;  Col = 4; New = "71÷5=" },
    @{ Row = 5;  Col = 5; New = "46÷9=" },

    @{ Row = 9;  Col = 1; New = "82÷8=" },
    @{ Row = 9;  Col = 2; New = "71÷4=" },
    @{ Row = 9;  Col = 3; New = "38÷5=" },
    @{ Row = 9;  Col = 4; New = "19÷8=" },
    @{ Row = 9;  Col = 5; New = "92÷2=" },

    @{ Row = 13; Col = 1; New = "16÷6=" },
    @{ Row = 13; Col = 2; New = "71÷7=" },
    @{ Row = 13; Col = 3; New = "46÷5=" },
    @{ Row = 13; Col = 4; New = "85÷5=" },
    @{ Row = 13; Col = 5; New = "15÷3=" },

    @{ Row = 17; Col = 1; New = "23÷7=" },
    @{ Row = 17; Col = 2; New = "10÷2=" },
    @{ Row = 17; Col = 3; New = "32÷3=" },
    @{ Row = 17; Col = 4; New = "64÷9=" },
    @{ Row = 17; Col = 5; New = "20÷2=" }
)

foreach ($change in $changes) {
    $cellRange = $table.Cell($change.Row, $change.Col).Range
    # Trim the trailing end-of-cell marker so we only overwrite the visible
    # text, leaving the cell's paragraph mark/formatting untouched.
    $cellRange.MoveEnd(12, -1) | Out-Null
    $cellRange.Text = $change.New
}
